$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.992.07"
$ws.Range("E2").Value = "  +1.24%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.205.35"
$ws.Range("E3").Value = "  +1.12%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.09%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.81"
$ws.Range("E5").Value = "  +3.91%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.16"
$ws.Range("E6").Value = "  +0.20%  "

# Row 7
$ws.Range("E7").Value = "  -0.06%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.202.38"
$ws.Range("E8").Value = "  +1.11%  "

# Row 9
$ws.Range("E9").Value = "  +1.45%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.162"
$ws.Range("E10").Value = "  -1.09%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.19"
$ws.Range("E11").Value = "  -0.09%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.510"
$ws.Range("E12").Value = "  +0.55%  "

# Row 13
$ws.Range("E13").Value = "  -0.91%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "38.20"
$ws.Range("E14").Value = "  -0.13%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.730.98"
$ws.Range("E15").Value = "  +1.19%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.610.45"
$ws.Range("E16").Value = "  +0.58%  "

# Row 17
$ws.Range("E17").Value = "  +1.10%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.200.49"
$ws.Range("E18").Value = "  +1.22%  "

# Row 19
$ws.Range("E19").Value = "  +0.20%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "513.80"
$ws.Range("E20").Value = "  +0.29%  "

# Row 21
$ws.Range("E21").Value = "  +6.00%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.733"
$ws.Range("E22").Value = "  -0.36%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "15.29"
$ws.Range("E23").Value = "  -4.59%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.97"
$ws.Range("E24").Value = "  +1.45%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.53"
$ws.Range("E25").Value = "  +0.68%  "

# Row 26
$ws.Range("E26").Value = "  +0.20%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.24"
$ws.Range("E27").Value = "  +1.02%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.00"
$ws.Range("E28").Value = "  +3.17%  "

# Row 29
$ws.Range("E29").Value = "  +1.54%  "

# Row 30
$ws.Range("E30").Value = "  +2.40%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "28.22"
$ws.Range("E31").Value = "  +0.56%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.24"
$ws.Range("E32").Value = "  +0.04%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.57"
$ws.Range("E33").Value = "  +3.72%  "

# Row 34
$ws.Range("E34").Value = "  -0.01%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.62"
$ws.Range("E35").Value = "  -0.83%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "55.27"
$ws.Range("E36").Value = "  -0.68%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0918"
$ws.Range("E37").Value = "  +4.18%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "485.00"
$ws.Range("E38").Value = "  +1.63%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0425"
$ws.Range("E39").Value = "  +0.45%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.04"
$ws.Range("E40").Value = "  -1.83%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.93"
$ws.Range("E41").Value = "  +3.35%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.020.43"
$ws.Range("E42").Value = "  -3.54%  "

# Row 43
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.118"
$ws.Range("E43").Value = "  -2.49%  "

# Row 44
$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.291"
$ws.Range("E44").Value = "  +0.16%  "

# Row 45
$ws.Range("E45").Value = "  -1.61%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₃0632"
$ws.Range("E46").Value = "  +7.65%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "29.07"
$ws.Range("E47").Value = "  -0.66%  "

# Row 49
$ws.Range("E49").Value = "  +0.42%  "

# Row 50
$ws.Range("E50").Value = "  -0.15%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "120.11"
$ws.Range("E51").Value = "  -2.85%  "
